$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.160.59'
$ws.Range('E2').Value = '  +0.22%  '

$ws.Range('D3').Value = '2.403.40'
$ws.Range('E3').Value = '  +5.28%  '

$ws.Range('E4').Value = '  -0.42%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '335.93'
$ws.Range('E5').Value = '  +9.01%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '105.34'
$ws.Range('E6').Value = '  -6.85%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.649'
$ws.Range('E7').Value = '  +2.70%  '

$ws.Range('E8').Value = '  -0.11%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.649'
$ws.Range('E9').Value = '  +5.80%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '42.27'
$ws.Range('E10').Value = '  -5.16%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0939'
$ws.Range('E11').Value = '  +1.38%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '8.73'
$ws.Range('E12').Value = '  -1.36%  '

$ws.Range('E13').Value = '  +1.02%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '17.36'
$ws.Range('E14').Value = '  +12.47%  '

$ws.Range('E15').Value = '  +1.95%  '

$ws.Range('D16').Value = '2.766.28'
$ws.Range('E16').Value = '  +5.40%  '

$ws.Range('D17').Value = '2.403.18'
$ws.Range('E17').Value = '  +5.28%  '

$ws.Range('D18').Value = '43.180.72'
$ws.Range('E18').Value = '  +0.36%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.71'
$ws.Range('E19').Value = '  +6.38%  '

$ws.Range('E20').Value = '  +0.99%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '3.94'
$ws.Range('E21').Value = '  +9.44%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '77.20'
$ws.Range('E22').Value = '  +2.53%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '275.53'
$ws.Range('E23').Value = '  +8.30%  '

$ws.Range('E24').Value = '  -2.13%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.87'
$ws.Range('E25').Value = '  +9.75%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '11.97'
$ws.Range('E26').Value = '  +1.98%  '

$ws.Range('E27').Value = '  +0.06%  '

$ws.Range('E28').Value = '  +4.98%  '

$ws.Range('E29').Value = '  -1.91%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '175.31'
$ws.Range('E30').Value = '  +0.29%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '37.09'
$ws.Range('E31').Value = '  -3.16%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.0940'
$ws.Range('E32').Value = '  +4.31%  '

$ws.Range('E33').Value = '  -1.01%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.08'
$ws.Range('E34').Value = '  +6.60%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.135'
$ws.Range('E35').Value = '  +5.02%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.89'
$ws.Range('E36').Value = '  -3.59%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.08'
$ws.Range('E37').Value = '  -3.36%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0367'
$ws.Range('E38').Value = '  -2.98%  '

$ws.Range('E39').Value = '  +3.67%  '

$ws.Range('E40').Value = '  +12.06%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.57'
$ws.Range('E41').Value = '  +12.89%  '

$ws.Range('E42').Value = '  +2.19%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '70.71'
$ws.Range('E43').Value = '  -2.44%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '121.85'
$ws.Range('E44').Value = '  +13.40%  '

$ws.Range('E45').Value = '  -0.03%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '92.22'
$ws.Range('E46').Value = '  +44.54%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '12.32'
$ws.Range('E47').Value = '  -2.66%  '

$ws.Range('E48').Value = '  -1.21%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.26'
$ws.Range('E49').Value = '  +5.29%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.506'
$ws.Range('E50').Value = '  +15.03%  '

$ws.Range('E51').Value = '  +1.50%  '
